## Adiciona link para video demo
## Slide 12, shape "TextBox 6": add a second line with the video demo link,
## hyperlinked to https://youtu.be/Q8tfYmYB7iw, and resize the textbox to fit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(2)

$tf = $shp.TextFrame
$tr = $tf.TextRange

# Append a new paragraph with the demo video URL, inheriting the existing run's
# character formatting (font, size, color).
$tr.Text = $tr.Text + "`r" + "https://youtu.be/Q8tfYmYB7iw"

$tr = $tf.TextRange

# The textbox is no longer centered now that it holds two lines of link text.
$tr.Paragraphs(1).ParagraphFormat.Alignment = [PpParagraphAlignment]::ppAlignLeft
$tr.Paragraphs(2).ParagraphFormat.Alignment = [PpParagraphAlignment]::ppAlignLeft

# Hyperlink the newly-added second line to the YouTube demo video.
$para2 = $tr.Paragraphs(2)
$para2.ActionSettings.Item(1).Hyperlink.Address = "https://youtu.be/Q8tfYmYB7iw"

# Grow the textbox height to fit the now two-line content (autofit).
$shp.Height = 36.35158
